$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.2393200005478144, 1.4252014926512548, -0.6163009051903825, 2.247312417022685)"
$ws.Range("C2").Value = "NCT(3.2099891108112537, 1.24663216574749, 2.6613372373170177, 4.506507713835729)"
$ws.Range("D2").Value = "JSU(-0.9929317501482251, 1.2475613393966303, 0.4685680304698515, 3.054054357077623)"
$ws.Range("E2").Value = "JSU(-1.2321670325439378, 1.0348201933930898, 2.927076423521182, 3.5562937944805615)"
